# Add a new slide after the existing slide 1, using the "Title Only"
# layout (ppLayoutTitleOnly = 11), matching the single title placeholder
# that appears in the target slide.
$p = $ppt.ActivePresentation

$newSlide = $p.Slides.Add(2, 11)

# Populate the title placeholder text.
$newSlide.Shapes.Item(1).TextFrame.TextRange.Text = "Try to learn"
